$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new values look numeric,
# so Excel stores them as exact text (matching the inlineStr values in the
# target workbook) instead of silently converting them to floating point numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the diff.
$ws.Range("D2").Value = "26.382.95"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "1.768.65"
$ws.Range("E3").Value = "  -2.55%  "
$ws.Range("D4").Value = "0.9975"
$ws.Range("E4").Value = "  -0.78%  "
$ws.Range("D5").Value = "0.9966"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "305.12"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").Value = "0.4272"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("E8").Value = "  +1.51%  "
$ws.Range("D9").Value = "0.07161"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "0.8473"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Value = "20.33"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.753.85"
$ws.Range("E12").Value = "  -4.56%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "6.422"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.235"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "0.06835"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "0.9974"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "78.74"
$ws.Range("E17").Value = "  -3.43%  "
$ws.Range("D18").Value = "0.000008670"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").Value = "0.9962"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "14.98"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "26.394.80"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").Value = "5.094"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "1.965.46"
$ws.Range("E24").Value = "  -4.68%  "
$ws.Range("D25").Value = "151.87"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").Value = "1.857"
$ws.Range("E26").Value = "  -5.74%  "
$ws.Range("D27").Value = "18.03"
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").Value = "113.71"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "1.805"
$ws.Range("E30").Value = "  +3.70%  "
$ws.Range("D31").Value = "0.08934"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").Value = "0.7279"
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("D33").Value = "1.122"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").Value = "4.317"
$ws.Range("E34").Value = "  -3.64%  "
$ws.Range("D35").Value = "0.9956"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "2.727"
$ws.Range("E36").Value = "  -6.65%  "
$ws.Range("D37").Value = "1.098"
$ws.Range("E37").Value = "  +2.03%  "
$ws.Range("D38").Value = "0.05148"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "0.01882"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").Value = "0.4927"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").Value = "0.1609"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D42").Value = "2.584"
$ws.Range("E42").Value = "  -6.48%  "
$ws.Range("D43").Value = "6.316"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "8.010"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").Value = "104.82"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("D47").Value = "0.9954"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "0.06191"
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("D49").Value = "0.4489"
$ws.Range("E49").Value = "  -2.78%  "
$ws.Range("D50").Value = "1.618"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").Value = "1.735"
$ws.Range("E51").Value = "  +2.38%  "
